# CryCompanywiseStockReport_1.xlsx — stock-qty/value recount.
# A handful of lines had their counted quantity (col F) and computed value
# (col G) corrected by -1 unit, a couple of duplicate-code item pairs got
# their two lot rows swapped (everything but SlNo/Description), and every
# Sub Total / Grand Total (col B) downstream was re-footed to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 49 (BHA-Ariel Liq 2L Pouch FL): qty 38 -> 37 ---
$ws.Range("F71").Value = 37
$ws.Range("G71").Value = 9279.23

# --- Row 78 (GLT-Mach 3 Turbo Cart 2s): qty 27 -> 26 ---
$ws.Range("F100").Value = 26
$ws.Range("G100").Value = 6397.82

# --- Sub Total: ---
$ws.Range("B114").Value = 289752.1

# --- Row 144 (EVE-Eveready AA Battery 1015): qty 2061 -> 2057 ---
$ws.Range("F186").Value = 2057
$ws.Range("G186").Value = 17381.65

# --- Row 145 (EVE-Eveready AAA Battery 1012): qty 1326 -> 1322 ---
$ws.Range("F187").Value = 1322
$ws.Range("G187").Value = 10562.78

# --- Sub Total: ---
$ws.Range("B189").Value = 35353.15

# --- Row 182 (HIM-BABY CARE GIFT PACK (WW)1): qty 25 -> 24 ---
$ws.Range("F234").Value = 24
$ws.Range("G234").Value = 7098

# --- Sub Total: ---
$ws.Range("B274").Value = 104571.73

# --- Rows 223/224 (HUL-Bru Inst Poly 50g): lot rows swapped ---
$ws.Range("B277").Value = 63565
$ws.Range("E277").Value = 109.19
$ws.Range("F277").Value = 60
$ws.Range("G277").Value = 6162.6
$ws.Range("B278").Value = 61610
$ws.Range("E278").Value = 122.71
$ws.Range("F278").Value = -58
$ws.Range("G278").Value = -5957.18

# --- Rows 240/242 (HUL-Kissan Pineapple Jam 500G): lot rows swapped ---
$ws.Range("B294").Value = 63531
$ws.Range("F294").Value = 80
$ws.Range("G294").Value = 11478.4
$ws.Range("B296").Value = 63571
$ws.Range("F296").Value = 8
$ws.Range("G296").Value = 1147.84

# --- Rows 245/246 (HUL-knorr schezwan 200g pch): lot rows swapped ---
$ws.Range("B299").Value = 63510
$ws.Range("E299").Value = 50.66
$ws.Range("F299").Value = 148
$ws.Range("G299").Value = 7050.72
$ws.Range("B300").Value = 55356
$ws.Range("E300").Value = 54.04
$ws.Range("F300").Value = -158
$ws.Range("G300").Value = -7527.12

# --- Rows 257/258 (HUL-lux advanced eventoned glow 4x100): lot rows swapped ---
$ws.Range("B311").Value = 63563
$ws.Range("E311").Value = 119.04
$ws.Range("F311").Value = 2
$ws.Range("G311").Value = 223.92
$ws.Range("B312").Value = 61605
$ws.Range("E312").Value = 133.78
$ws.Range("F312").Value = -13
$ws.Range("G312").Value = -1455.48

# --- Row 287 (HUL-Hlx Bib 500g): qty 5 -> 4 ---
$ws.Range("F343").Value = 4
$ws.Range("G343").Value = 646.92

# --- Sub Total: ---
$ws.Range("B346").Value = 10804.08

# --- Row 304 (JYOTHY-Fa Mens Deo Mix 150 ml): qty 35 -> 34 ---
$ws.Range("F366").Value = 34
$ws.Range("G366").Value = 4668.88

# --- Sub Total: ---
$ws.Range("B395").Value = 264045.53

# --- Rows 350/351 (KUS-Floor Wiper): lot rows swapped ---
$ws.Range("B420").Value = 58047
$ws.Range("D420").Value = 105.54
$ws.Range("E420").Value = 126.1
$ws.Range("F420").Value = 43
$ws.Range("G420").Value = 4538.22
$ws.Range("B421").Value = 47097
$ws.Range("D421").Value = 112.28
$ws.Range("E421").Value = 134.16
$ws.Range("F421").Value = 15
$ws.Range("G421").Value = 1684.2

# --- Rows 387/388 (CRE-Bourbon 100gm): lot rows swapped ---
$ws.Range("B465").Value = 65069
$ws.Range("E465").Value = 14.3
$ws.Range("F465").Value = 2
$ws.Range("G465").Value = 26.9
$ws.Range("B466").Value = 53757
$ws.Range("E466").Value = 16.08
$ws.Range("F466").Value = -159
$ws.Range("G466").Value = -2138.55

# --- Rows 398/399 (CRE-Cremica Golden Bytes Rich Butter 200Gm): lot rows swapped ---
$ws.Range("B476").Value = 64922
$ws.Range("E476").Value = 20.98
$ws.Range("F476").Value = 136
$ws.Range("G476").Value = 2683.28
$ws.Range("B477").Value = 45706
$ws.Range("E477").Value = 23.58
$ws.Range("F477").Value = -202
$ws.Range("G477").Value = -3985.46

# --- Rows 401/402 (CRE-Cremica Honey Oatmeal Cookies 50 +25 Gm): lot rows swapped ---
$ws.Range("B479").Value = 45718
$ws.Range("E479").Value = 19.38
$ws.Range("F479").Value = -294
$ws.Range("G479").Value = -4768.68
$ws.Range("B480").Value = 64927
$ws.Range("E480").Value = 17.26
$ws.Range("F480").Value = 217
$ws.Range("G480").Value = 3519.74

# --- Rows 412/413 (CRE-Kaju khz cookies 100 gm): lot rows swapped ---
$ws.Range("B490").Value = 65067
$ws.Range("E490").Value = 15.65
$ws.Range("F490").Value = 252
$ws.Range("G490").Value = 3711.96
$ws.Range("B491").Value = 53595
$ws.Range("E491").Value = 17.61
$ws.Range("F491").Value = -335
$ws.Range("G491").Value = -4934.55

# --- Rows 470/471 (PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)): lot rows swapped ---
$ws.Range("B564").Value = 53319
$ws.Range("E564").Value = 310.64
$ws.Range("F564").Value = -6
$ws.Range("G564").Value = -1643.52
$ws.Range("B565").Value = 64810
$ws.Range("E565").Value = 291.22
$ws.Range("F565").Value = 6
$ws.Range("G565").Value = 1643.52

# --- Rows 498/499 (Rasna Nagpur Orange (32 Glass)): lot rows swapped ---
$ws.Range("B596").Value = 64830
$ws.Range("E596").Value = 34.9
$ws.Range("F596").Value = 113
$ws.Range("G596").Value = 3709.79
$ws.Range("B597").Value = 60022
$ws.Range("E597").Value = 37.22
$ws.Range("F597").Value = -113
$ws.Range("G597").Value = -3709.79

# --- Row 519 (Sarathi-Laya 125 Gms Zipper Pouch): qty 12 -> 11 ---
$ws.Range("F627").Value = 11
$ws.Range("G627").Value = 474.98

# --- Row 520 (Sarathi-Sii-Manpasad 125 Grams Zipper): qty 19 -> 18 ---
$ws.Range("F628").Value = 18
$ws.Range("G628").Value = 777.24

# --- Sub Total: ---
$ws.Range("B629").Value = 1268.08

# --- Row 585 (TCP-Urad Dal 1 kg): qty 365 -> 364 ---
$ws.Range("F715").Value = 364
$ws.Range("G715").Value = 43938.44

# --- Sub Total: ---
$ws.Range("B716").Value = 190532.32

# --- Rows 600/601 (Shankys Tip Top Hing Jeera Peanut/ Salted Peanut 200 Gm): lot rows swapped ---
$ws.Range("B732").Value = 65079
$ws.Range("F732").Value = 21
$ws.Range("G732").Value = 858.27
$ws.Range("B733").Value = 65362
$ws.Range("F733").Value = 69
$ws.Range("G733").Value = 2820.03

# --- Row 610 (Tip Top Sooji 1 Kg): qty 130 -> 129 ---
$ws.Range("F742").Value = 129
$ws.Range("G742").Value = 7351.71

# --- Sub Total: ---
$ws.Range("B743").Value = 103964.18

# --- Row 629 (VVD Pure Drop Cold Pressed Gingelly Oil Pouch 500Ml): qty 517 -> 515 ---
$ws.Range("F771").Value = 515
$ws.Range("G771").Value = 74494.75

# --- Row 632 (VVD Veda Pancha Deepam Oil Pouch 900Ml): qty 241 -> 240 ---
$ws.Range("F774").Value = 240
$ws.Range("G774").Value = 30859.2

# --- Sub Total: ---
$ws.Range("B775").Value = 863150.11

# --- Sub Total: / Grand Total: ---
$ws.Range("B793").Value = 3330749.96
$ws.Range("B794").Value = 3330749.96
